$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.26'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.06'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.407'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06005'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.390'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8084'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9299'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1419'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07438'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03370'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03033'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09358'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.940'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001589'

$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0005942'
$ws.Range("E17").Value = '16OneONE'

$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005467'
$ws.Range("E18").Value = '17TigerCashTCH'

$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004152'
$ws.Range("E19").Value = '18HotbitTokenHTB'

$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009881'
$ws.Range("E20").Value = '19BitKanKAN'

$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.00007103'
$ws.Range("E21").Value = '20NitroExNTX'

$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.663'
$ws.Range("E22").Value = '21LEOLEO'

$ws.Range("B23").Value = 'KuCoinToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.432'
$ws.Range("E23").Value = '22KuCoinTokenKCS'

$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.186'
$ws.Range("E24").Value = '23BTSETokenBTSE'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03968'

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1075'
$ws.Range("E41").Value = '40BKEXTokenBKK'

$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002711'
$ws.Range("E42").Value = '41CEJICEJI'

$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003034'
$ws.Range("E43").Value = '42KickTokenKICK'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.006202'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005196'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0005802'
